# Auto-generated edit script: updates computed profit/price columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR item-flipping sheets to reflect
# refreshed Universalis market data from the scheduled runner.
$wb = $excel.ActiveWorkbook

# ALC sheet, row 104
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(104, 8).Value = 284.5
$ws.Cells.Item(104, 9).Value = 284.5
$ws.Cells.Item(104, 11).Value = 853.5
$ws.Cells.Item(104, 13).Value = 893.5

# ALC sheet, row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 3940
$ws.Cells.Item(113, 9).Value = 3800
$ws.Cells.Item(113, 10).Value = 4500
$ws.Cells.Item(113, 11).Value = 3800
$ws.Cells.Item(113, 12).Value = 4500
$ws.Cells.Item(113, 13).Value = -546
$ws.Cells.Item(113, 14).Value = -11008

# ALC sheet, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 9256.825999999999
$ws.Cells.Item(132, 9).Value = 11528.444
$ws.Cells.Item(132, 11).Value = 34585.33199999999
$ws.Cells.Item(132, 13).Value = -32055.33199999999

# ALC sheet, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 4033.7908
$ws.Cells.Item(138, 9).Value = 3309.6553
$ws.Cells.Item(138, 10).Value = 5533.7856
$ws.Cells.Item(138, 11).Value = 9928.965899999999
$ws.Cells.Item(138, 12).Value = 16601.3568
$ws.Cells.Item(138, 13).Value = -4788.965899999999
$ws.Cells.Item(138, 14).Value = -26881.3568

# ARM sheet, row 10
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 14).ClearContents()

# ARM sheet, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 250389
$ws.Cells.Item(32, 9).Value = 254148.72
$ws.Cells.Item(32, 11).Value = 254148.72
$ws.Cells.Item(32, 13).Value = -253861.72

# ARM sheet, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 6153.923
$ws.Cells.Item(61, 9).Value = 5166.8335
$ws.Cells.Item(61, 10).Value = 7000
$ws.Cells.Item(61, 11).Value = 5166.8335
$ws.Cells.Item(61, 12).Value = 7000
$ws.Cells.Item(61, 13).Value = -4954.8335
$ws.Cells.Item(61, 14).Value = -7424

# ARM sheet, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 6900.5654
$ws.Cells.Item(132, 9).Value = 4888.478
$ws.Cells.Item(132, 11).Value = 14665.434
$ws.Cells.Item(132, 13).Value = -12135.434

# ARM sheet, row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(133, 14).ClearContents()

# ARM sheet, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 6153.923
$ws.Cells.Item(136, 9).Value = 5166.8335
$ws.Cells.Item(136, 10).Value = 7000
$ws.Cells.Item(136, 11).Value = 15500.5005
$ws.Cells.Item(136, 12).Value = 21000
$ws.Cells.Item(136, 13).Value = -12950.5005
$ws.Cells.Item(136, 14).Value = -26100

# BSM sheet, row 26
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 45490.332
$ws.Cells.Item(26, 10).Value = 65000
$ws.Cells.Item(26, 12).Value = 65000
$ws.Cells.Item(26, 14).Value = -65584

# BSM sheet, row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 6522.25
$ws.Cells.Item(105, 9).Value = 2016.25
$ws.Cells.Item(105, 11).Value = 2016.25
$ws.Cells.Item(105, 13).Value = -269.25

# BSM sheet, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 6388.5557
$ws.Cells.Item(134, 9).Value = 8874.5
$ws.Cells.Item(134, 10).Value = 4399.8
$ws.Cells.Item(134, 11).Value = 26623.5
$ws.Cells.Item(134, 12).Value = 13199.4
$ws.Cells.Item(134, 13).Value = -24088.5
$ws.Cells.Item(134, 14).Value = -18269.4

# BSM sheet, row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(135, 8).Value = 69459.5
$ws.Cells.Item(135, 10).Value = 69459.5
$ws.Cells.Item(135, 12).Value = 69459.5
$ws.Cells.Item(135, 14).Value = -79599.5

# CRP sheet, row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 104489.45
$ws.Cells.Item(16, 9).Value = 26000
$ws.Cells.Item(16, 11).Value = 26000
$ws.Cells.Item(16, 13).Value = -25713

# CRP sheet, row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 17489
$ws.Cells.Item(99, 9).Value = 17489
$ws.Cells.Item(99, 11).Value = 17489
$ws.Cells.Item(99, 13).Value = -15991

# CRP sheet, row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 104489.45
$ws.Cells.Item(113, 9).Value = 26000
$ws.Cells.Item(113, 11).Value = 26000
$ws.Cells.Item(113, 13).Value = -23830

# CRP sheet, row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 17489
$ws.Cells.Item(126, 9).Value = 17489
$ws.Cells.Item(126, 11).Value = 52467
$ws.Cells.Item(126, 13).Value = -49997

# CRP sheet, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2443.0605
$ws.Cells.Item(132, 9).Value = 2471.8386
$ws.Cells.Item(132, 10).Value = 1997
$ws.Cells.Item(132, 11).Value = 7415.5158
$ws.Cells.Item(132, 12).Value = 5991
$ws.Cells.Item(132, 13).Value = -4885.5158
$ws.Cells.Item(132, 14).Value = -11051

# CRP sheet, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 3367.5625
$ws.Cells.Item(134, 9).Value = 2990.4
$ws.Cells.Item(134, 11).Value = 8971.200000000001
$ws.Cells.Item(134, 13).Value = -6436.200000000001

# CUL sheet, row 93
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(93, 8).Value = 5000
$ws.Cells.Item(93, 10).Value = 5000
$ws.Cells.Item(93, 12).Value = 15000
$ws.Cells.Item(93, 14).Value = -18744

# CUL sheet, row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 859.3182
$ws.Cells.Item(113, 10).Value = 910.5263
$ws.Cells.Item(113, 12).Value = 2731.5789
$ws.Cells.Item(113, 14).Value = -7071.5789

# CUL sheet, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 769017.1
$ws.Cells.Item(122, 9).Value = 5376623.5
$ws.Cells.Item(122, 10).Value = 1082.7222
$ws.Cells.Item(122, 11).Value = 48389611.5
$ws.Cells.Item(122, 12).Value = 9744.4998
$ws.Cells.Item(122, 13).Value = -48387161.5
$ws.Cells.Item(122, 14).Value = -14644.4998

# CUL sheet, row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 502820.84
$ws.Cells.Item(129, 9).Value = 1431716.9
$ws.Cells.Item(129, 10).Value = 2646.077
$ws.Cells.Item(129, 11).Value = 4295150.699999999
$ws.Cells.Item(129, 12).Value = 7938.231000000001
$ws.Cells.Item(129, 13).Value = -4290150.699999999
$ws.Cells.Item(129, 14).Value = -17938.231

# CUL sheet, row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 7474.227
$ws.Cells.Item(139, 9).Value = 6586.8
$ws.Cells.Item(139, 10).Value = 7735.2354
$ws.Cells.Item(139, 11).Value = 19760.4
$ws.Cells.Item(139, 12).Value = 23205.7062
$ws.Cells.Item(139, 13).Value = -14620.4
$ws.Cells.Item(139, 14).Value = -33485.7062

# CUL sheet, row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 1974.25
$ws.Cells.Item(140, 9).Value = 1974.25
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 5922.75
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 13).Value = -742.75
$ws.Cells.Item(140, 14).ClearContents()

# GSM sheet, row 26
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 14).ClearContents()

# GSM sheet, row 50
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 14).ClearContents()

# GSM sheet, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2563.6667
$ws.Cells.Item(102, 9).Value = 2563.6667
$ws.Cells.Item(102, 11).Value = 2563.6667
$ws.Cells.Item(102, 13).Value = -941.6667000000002

# GSM sheet, row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 376.93332
$ws.Cells.Item(107, 9).Value = 383.41666
$ws.Cells.Item(107, 10).Value = 351
$ws.Cells.Item(107, 11).Value = 383.41666
$ws.Cells.Item(107, 12).Value = 351
$ws.Cells.Item(107, 13).Value = 1536.58334
$ws.Cells.Item(107, 14).Value = -4191

# GSM sheet, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3798.45
$ws.Cells.Item(122, 9).Value = 4116.25
$ws.Cells.Item(122, 11).Value = 12348.75
$ws.Cells.Item(122, 13).Value = -9898.75

# LTW sheet, row 53
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 13).ClearContents()

# LTW sheet, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4046.25
$ws.Cells.Item(132, 9).Value = 3231
$ws.Cells.Item(132, 11).Value = 9693
$ws.Cells.Item(132, 13).Value = -7163

# WVR sheet, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 32451.082
$ws.Cells.Item(122, 9).Value = 2267.5386
$ws.Cells.Item(122, 11).Value = 6802.6158
$ws.Cells.Item(122, 13).Value = -4352.6158

# WVR sheet, row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2921
$ws.Cells.Item(126, 9).Value = 2666.6667
$ws.Cells.Item(126, 10).Value = 3302.5
$ws.Cells.Item(126, 11).Value = 8000.000100000001
$ws.Cells.Item(126, 12).Value = 9907.5
$ws.Cells.Item(126, 13).Value = -5530.000100000001
$ws.Cells.Item(126, 14).Value = -14847.5

# WVR sheet, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2526.7778
$ws.Cells.Item(132, 9).Value = 1541.5
$ws.Cells.Item(132, 11).Value = 4624.5
$ws.Cells.Item(132, 13).Value = -2094.5
